$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column M (13) — this shifts M..V to N..W,
# preserving styles/values, and updates the sheet dimension automatically.
$ws.Columns.Item(13).Insert()

# Header for the newly inserted column.
$ws.Range("M1").Value = "fuel_remaining(liters)"

# Fuel remaining (liters) values for the data rows (the blank/divider rows
# 4, 7, 11, 14, 17-20 stay empty).
$ws.Range("M2").Value = 244
$ws.Range("M3").Value = 78
$ws.Range("M5").Value = 815
$ws.Range("M6").Value = 231
$ws.Range("M8").Value = 713
$ws.Range("M9").Value = 636
$ws.Range("M10").Value = 154
$ws.Range("M12").Value = 1079
$ws.Range("M13").Value = 757
$ws.Range("M15").Value = 1234
$ws.Range("M16").Value = 78

# Update the "description of work" text (now column W) for row 8 to mention
# bringing a fuel barrel.
$ws.Range("W8").Value = "long day. bring 1 x 200 l fuel barrel with. new AWS install at Neem and visit Humboldt after NEEM landing and ground stop. no crane?"
